$d = $word.ActiveDocument

# Helper: scoped literal find & replace within a given Range, so that
# duplicate text elsewhere in the document is never touched.
function Replace-InRange($range, [string]$old, [string]$new) {
    $ok = $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Find failed for: $old"
    }
}

# ------------------------------------------------------------------
# 1) "Objetivos" section (paragraphs 6 & 7, 1-based) takes the text
#    that used to belong to "Programa resumido" (paragraphs 11 & 12).
# ------------------------------------------------------------------
Replace-InRange $d.Paragraphs(6).Range `
    'Apresentar aos alunos os fundamentos do planejamento e da gestão estratégica nas organizações, capacitando-os quanto as metodologias existentes, suas etapas e implicações para os resultados organizacionais.' `
    'Administração e processo estratégico; Planejamento estratégico, tático e operacional; Diretrizes organizacionais: missão, visão e objetivos; Formulação estratégica; Implantação de estratégia; Controle estratégico; Planejamento de unidades organizacionais; Administração estratégica aplicada.'

Replace-InRange $d.Paragraphs(7).Range `
    'Introduce students to the fundamentals of planning and strategic management in organizations, training them as to the existing methodologies, their stages and implications for organizational results.' `
    'Management and strategic process; Strategic, tactical and operational planning; Organizational guidelines: mission, vision and objectives; Strategic formulation; Strategy implementation; Strategic control; Planning of organizational units; Strategic management applied.'

# ------------------------------------------------------------------
# 2) "Docente(s) Responsável(eis)" list bullet (paragraph 9) takes the
#    text that used to be the "Objetivos" paragraph.
# ------------------------------------------------------------------
Replace-InRange $d.Paragraphs(9).Range `
    '11079086 - Herlandí de Souza Andrade' `
    'Apresentar aos alunos os fundamentos do planejamento e da gestão estratégica nas organizações, capacitando-os quanto as metodologias existentes, suas etapas e implicações para os resultados organizacionais.'

# ------------------------------------------------------------------
# 3) "Programa resumido" body (paragraphs 11 & 12) takes the text that
#    used to belong to "Programa" (paragraphs 14 & 15).
# ------------------------------------------------------------------
Replace-InRange $d.Paragraphs(11).Range `
    'Administração e processo estratégico; Planejamento estratégico, tático e operacional; Diretrizes organizacionais: missão, visão e objetivos; Formulação estratégica; Implantação de estratégia; Controle estratégico; Planejamento de unidades organizacionais; Administração estratégica aplicada.' `
    '1. Motivações e Desafios para a estratégia; 2. Conceitos Básicos de Estratégia; 3. Gestão Estratégica; 4. Transformação Estratégica; 5. Análise do Ambiente Externo; 6. Análise da Turbulência e da Vulnerabilidade; 7. Análise do Ambiente Interno; 8. Representação do Portifólio; 9. Estratégia de Balanceamento do Portifólio; 10. Formulação das Estratégias; 11. Capacitação Estratégica; 12. O Plano Estratégico; 13. Metodologia de Planejamento Estratégico; 14. Workshop de Planejamento Estratégico; 15. Implantação da Gestão Estratégica.'

Replace-InRange $d.Paragraphs(12).Range `
    'Management and strategic process; Strategic, tactical and operational planning; Organizational guidelines: mission, vision and objectives; Strategic formulation; Strategy implementation; Strategic control; Planning of organizational units; Strategic management applied.' `
    'Introduce students to the fundamentals of planning and strategic management in organizations, training them as to the existing methodologies, their stages and implications for organizational results.'

# ------------------------------------------------------------------
# 4) "Programa" body (paragraph 14) takes the text that used to be the
#    "Método:" value inside "Avaliação".
# ------------------------------------------------------------------
Replace-InRange $d.Paragraphs(14).Range `
    '1. Motivações e Desafios para a estratégia; 2. Conceitos Básicos de Estratégia; 3. Gestão Estratégica; 4. Transformação Estratégica; 5. Análise do Ambiente Externo; 6. Análise da Turbulência e da Vulnerabilidade; 7. Análise do Ambiente Interno; 8. Representação do Portifólio; 9. Estratégia de Balanceamento do Portifólio; 10. Formulação das Estratégias; 11. Capacitação Estratégica; 12. O Plano Estratégico; 13. Metodologia de Planejamento Estratégico; 14. Workshop de Planejamento Estratégico; 15. Implantação da Gestão Estratégica.' `
    'Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras.'

# ------------------------------------------------------------------
# 5) "Avaliação" list bullet (paragraph 17): each value run shifts to
#    the next label, and the bibliography list gets folded into the
#    last ("Norma de recuperação:") run, joined with manual line
#    breaks (character 11, the same code Word uses for <w:br/>)
#    two-by-two, exactly like the original "Bibliografia" paragraph.
# ------------------------------------------------------------------
$avaliacaoRange = $d.Paragraphs(17).Range

Replace-InRange $avaliacaoRange `
    'Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras.' `
    'Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas'

Replace-InRange $avaliacaoRange `
    'Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas' `
    'NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.'

$vt = [char]11
$bibliografia = @(
    'OLIVEIRA, D. P. R. Planejamento Estratégico: Conceitos, Metodologia e Práticas. 34 ed. São Paulo: Atlas, 2018.',
    'CHIAVENATO, I; SAPIRO, A. Planejamento Estratégico: Fundamentos e Aplicações. 3 ed. Rio de Janeiro: Campus, 2015.',
    'COSTA, E. A. Gestão Estratégica: da empresa que temos para a empresa que queremos. 2 ed. Santo André: Saraiva, 2012.',
    'LOBATO, D. M. Estratégia de Empresas. Rio de Janeiro: FGV, 2009.',
    'HITT, M A. et al. Administração Estratégica. São Paulo: Pioneira Thomson Learning, 2007.',
    'GHEMAWAT, P. A Estratégia e o cenário de Negócios. Porto Alegre: Bookman, 2007.',
    'MINTZBERG, H. et al. O Processo da Estratégia. São Paulo: Bookman, 2006.',
    'HAMEL, G., PRAHALAD, C.K. Competindo pelo futuro. Rio de Janeiro: Campus, 2005.',
    'PORTER, M. Estratégia Competitiva. Rio de janeiro: Campus, 2005.',
    'KAPLAN, R. S. Mapas Estratégicos: Balanced Scorecard. Rio de Janeiro: Elsevier, 2004.'
)
$bibliografiaText = [string]::Join($vt + $vt, $bibliografia)
$normaOld = 'NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.'
$normaNew = $normaOld + $vt + $vt + $bibliografiaText

Replace-InRange $avaliacaoRange $normaOld $normaNew

# ------------------------------------------------------------------
# 6) The old "Bibliografia" body paragraph (paragraph 19) now only
#    holds the professor entry that used to be the "Docente(s)" bullet.
#    Its old content spanned many runs/line-breaks, so replace the
#    paragraph's whole range (excluding the trailing paragraph mark)
#    rather than a single Find match.
# ------------------------------------------------------------------
$bibParaRange = $d.Paragraphs(19).Range
$bibParaRange.MoveEnd(1, -1) | Out-Null
$bibParaRange.Text = '11079086 - Herlandí de Souza Andrade'
